# Pre treatment phase outcome measurements added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row labels that changed text
$ws.Range("A11").Value = "Task-Monitor Scale"
$ws.Range("A12").Value = "Organization of Materials Scale"

# Add the new "Pre Experimental Phase" (column C) data for rows 2-14
$ws.Range("C2").Value = 62
$ws.Range("C3").Value = 65
$ws.Range("C4").Value = 64
$ws.Range("C5").Value = 68
$ws.Range("C6").Value = 62
$ws.Range("C7").Value = 66
$ws.Range("C8").Value = 70
$ws.Range("C9").Value = 70
$ws.Range("C10").Value = 58
$ws.Range("C11").Value = 61
$ws.Range("C12").Value = 49
$ws.Range("C13").Value = 64
$ws.Range("C14").Value = 65

# Widen column A to account for the longer label text
$ws.Columns.Item(1).ColumnWidth = 28.83203125

# Update selection to match final state
$ws.Range("B26").Select() | Out-Null
